# "edições para trabalhar no caos"
# Turn the numeric header row (1..9) in B1:J1 into text labels "P1".."P9",
# matching the look of the ID header in A1 (bold/plain text header row),
# and apply a (black) font/style to the whole header row A1:J1.
# Finally, reproduce the new header-row selection (A1:J1) left behind by
# the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric sequence in the header row with P1..P9 labels.
$headers = @("P1", "P2", "P3", "P4", "P5", "P6", "P7", "P8", "P9")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Style the whole header row (A1:J1) with the new font (explicit black
# RGB color rather than the theme color used before).
$headerRow = $ws.Range("A1:J1")
$headerRow.Font.Color = 0

# Leave the header row selected, like after the edit was made.
$ws.Range("A1:J1").Select() | Out-Null
